$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.107.65'
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = '  -2.27%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.821.04'
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = '  -1.49%  '
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -1.26%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '310.63'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -3.07%  '
$ws.Cells.Item(6, 5).Value = '  -1.14%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4225'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -1.92%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3660'
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -2.08%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07221'
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -2.02%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8468'
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -3.61%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.88'
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -3.48%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.828.10'
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -1.23%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.646'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -1.23%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07082'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -0.70%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.277'
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -3.22%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '89.33'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +1.16%  '
$ws.Cells.Item(17, 5).Value = '  -1.30%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000008825'
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -1.95%  '
$ws.Cells.Item(20, 2).Value = 'BitDAO'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5078'
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -1.98%  '
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '14.97'
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -3.25%  '
$ws.Cells.Item(22, 2).Value = 'WrappedBTC'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.164.35'
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -2.09%  '
$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.103'
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -2.56%  '
$ws.Cells.Item(24, 2).Value = 'Cosmos'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.82'
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  -2.58%  '
$ws.Cells.Item(25, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.051.25'
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -1.81%  '
$ws.Cells.Item(26, 2).Value = 'Toncoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.973'
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -1.77%  '
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '151.87'
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -2.62%  '
$ws.Cells.Item(28, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.241'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +4.62%  '
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.31'
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  -1.71%  '
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.209'
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -3.54%  '
$ws.Cells.Item(31, 2).Value = 'BitcoinCash'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '116.03'
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -2.65%  '
$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.08792'
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -1.94%  '
$ws.Cells.Item(33, 2).Value = 'ARBITRUM'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.178'
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -4.44%  '
$ws.Cells.Item(34, 5).Value = '  +1.41%  '
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7411'
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -4.72%  '
$ws.Cells.Item(36, 2).Value = 'Filecoin'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.420'
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -3.30%  '
$ws.Cells.Item(37, 2).Value = 'Frax'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -1.21%  '
$ws.Cells.Item(38, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.091'
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -4.03%  '
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.01962'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  -0.55%  '
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05236'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -2.24%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.273'
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +0.05%  '
$ws.Cells.Item(42, 2).Value = 'MXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.870'
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -0.66%  '
$ws.Cells.Item(43, 2).Value = 'Algorand'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1688'
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +0.00%  '
$ws.Cells.Item(44, 2).Value = 'TheSandbox'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5020'
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -2.61%  '
$ws.Cells.Item(45, 2).Value = 'Aptos'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.576'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -3.08%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.60'
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -0.59%  '
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4740'
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -0.11%  '
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '106.14'
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -3.09%  '
$ws.Cells.Item(49, 2).Value = 'PaxDollar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -1.30%  '
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06368'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -1.87%  '
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.655'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -2.30%  '
